$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a Price-column (D) cell as text first (NumberFormat "@") so the
# dotted/decimal-like strings below are not auto-converted to numbers by
# Excel (e.g. "1.000" -> 1, or "0.00001100" -> 1.1E-05). Each cell is set
# individually - multi-area Range(...) NumberFormat assignment is not
# reliable in this host.
function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
}

# Row 2 - Bitcoin
Set-TextValue "D2" "28.178.33"
$ws.Range("E2").Value = "  +1.11%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.807.88"
$ws.Range("E3").Value = "  +3.88%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.51%  "

# Row 5 - BNB
Set-TextValue "D5" "335.34"
$ws.Range("E5").Value = "  -0.08%  "

# Row 6 - USDC
Set-TextValue "D6" "0.9993"
$ws.Range("E6").Value = "  -0.32%  "

# Row 7 - XRP
Set-TextValue "D7" "0.4702"
$ws.Range("E7").Value = "  +25.09%  "

# Row 8 - Cardano
Set-TextValue "D8" "0.3721"
$ws.Range("E8").Value = "  +11.51%  "

# Row 9 - OKB
Set-TextValue "D9" "45.35"
$ws.Range("E9").Value = "  +1.99%  "

# Row 10 - Dogecoin
Set-TextValue "D10" "0.07701"
$ws.Range("E10").Value = "  +7.30%  "

# Row 11 - Polygon
Set-TextValue "D11" "1.157"
$ws.Range("E11").Value = "  +4.52%  "

# Row 12 - Solana
Set-TextValue "D12" "22.63"
$ws.Range("E12").Value = "  +1.70%  "

# Row 13 - BinanceUSD
Set-TextValue "D13" "1.000"
$ws.Range("E13").Value = "  -0.42%  "

# Row 14 - Polkadot
Set-TextValue "D14" "6.391"
$ws.Range("E14").Value = "  +4.48%  "

# Row 15 - Chainlink
Set-TextValue "D15" "7.409"
$ws.Range("E15").Value = "  +4.63%  "

# Row 16 - WrappedEther
Set-TextValue "D16" "1.800.64"
$ws.Range("E16").Value = "  +3.27%  "

# Row 17 - ShibaInu
Set-TextValue "D17" "0.00001100"
$ws.Range("E17").Value = "  +4.66%  "

# Row 18 - TRON
Set-TextValue "D18" "0.06733"
$ws.Range("E18").Value = "  +2.79%  "

# Row 19 - Litecoin
Set-TextValue "D19" "82.83"
$ws.Range("E19").Value = "  +5.26%  "

# Row 20 - Dai
Set-TextValue "D20" "0.9997"
$ws.Range("E20").Value = "  -0.34%  "

# Row 21 - Avalanche
$ws.Range("E21").Value = "  +4.78%  "

# Row 22 - Uniswap
Set-TextValue "D22" "6.446"
$ws.Range("E22").Value = "  +3.79%  "

# Row 23 - WrappedBTC
Set-TextValue "D23" "28.177.74"
$ws.Range("E23").Value = "  +1.01%  "

# Row 24 - Cosmos
$ws.Range("E24").Value = "  +3.07%  "

# Row 25 - Toncoin
Set-TextValue "D25" "2.409"
$ws.Range("E25").Value = "  +0.53%  "

# Row 26 - EthereumClassic
Set-TextValue "D26" "21.03"
$ws.Range("E26").Value = "  +6.90%  "

# Row 27 - LidoDAOToken
Set-TextValue "D27" "2.417"
$ws.Range("E27").Value = "  +5.14%  "

# Row 28 - Monero
Set-TextValue "D28" "153.30"
$ws.Range("E28").Value = "  +0.72%  "

# Row 29 - WrappedliquidstakedEther2.0
Set-TextValue "D29" "2.007.41"
$ws.Range("E29").Value = "  +3.29%  "

# Row 30 - BitcoinCash
Set-TextValue "D30" "134.64"
$ws.Range("E30").Value = "  +2.65%  "

# Row 31 - ImmutableX
Set-TextValue "D31" "1.275"
$ws.Range("E31").Value = "  +1.93%  "

# Row 32 - HuobiToken
Set-TextValue "D32" "4.044"
$ws.Range("E32").Value = "  +0.40%  "

# Row 33 - Stellar
$ws.Range("E33").Value = "  +11.16%  "

# Row 34 - Filecoin
Set-TextValue "D34" "5.953"
$ws.Range("E34").Value = "  +3.72%  "

# Row 35 - Algorand
Set-TextValue "D35" "0.2248"
$ws.Range("E35").Value = "  +6.84%  "

# Row 36 - Aptos
Set-TextValue "D36" "12.32"
$ws.Range("E36").Value = "  +1.31%  "

# Row 37 - was VeChain, now Hedera (rows 37/38 swapped order)
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D37" "0.06428"
$ws.Range("E37").Value = "  +3.92%  "

# Row 38 - was Hedera, now VeChain
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D38" "0.02382"
$ws.Range("E38").Value = "  +3.18%  "

# Row 39 - TheSandbox
Set-TextValue "D39" "0.6762"
$ws.Range("E39").Value = "  +2.15%  "

# Row 40 - InternetComputer(DFINITY)
Set-TextValue "D40" "5.283"
$ws.Range("E40").Value = "  +2.90%  "

# Row 41 - WEMIXTOKEN
Set-TextValue "D41" "1.512"
$ws.Range("E41").Value = "  +4.90%  "

# Row 42 - TrustWalletToken
Set-TextValue "D42" "1.235"
$ws.Range("E42").Value = "  +2.54%  "

# Row 43 - FraxShare
Set-TextValue "D43" "8.156"
$ws.Range("E43").Value = "  +3.33%  "

# Row 44 - EnergySwap
$ws.Range("E44").Value = "  +3.35%  "

# Row 45 - Frax
Set-TextValue "D45" "0.9989"

# Row 46 - Decentraland
Set-TextValue "D46" "0.6211"
$ws.Range("E46").Value = "  +3.52%  "

# Row 47 - PancakeSwap
Set-TextValue "D47" "3.840"
$ws.Range("E47").Value = "  +0.80%  "

# Row 48 - Quant
Set-TextValue "D48" "130.23"
$ws.Range("E48").Value = "  +2.71%  "

# Row 49 - NEARProtocol
Set-TextValue "D49" "2.077"
$ws.Range("E49").Value = "  +3.69%  "

# Row 50 - EOS
$ws.Range("E50").Value = "  +2.96%  "

# Row 51 - Cronos
Set-TextValue "D51" "0.07158"
$ws.Range("E51").Value = "  +1.06%  "
